# Apply the recorded change: a cyclic rotation of the data held in rows 20-22
# (new row20 <= old row21, new row21 <= old row22, new row22 <= old row20),
# expressed here as the concrete per-cell value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preventing Excel's "General"
# number-format auto-detection from turning date-looking strings (e.g.
# "2023-08-18") into date serials. We temporarily force Text format, assign
# the value, then restore the cell's original number format.
function Set-TextValue($range, $value) {
    $originalFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $originalFormat
}

# ---- Row 20 ----
$ws.Range("A20").Value = 111560002
$ws.Range("B20").Value = 94134
Set-TextValue $ws.Range("D20") "NT"
$ws.Range("E20").Value = 53
Set-TextValue $ws.Range("F20") "Vedtrappmossa"
Set-TextValue $ws.Range("G20") "Crossocalyx hellerianus"
Set-TextValue $ws.Range("H20") "(Nees ex Lindenb.) Meyl."
$ws.Range("Q20").Value = 523901.5626756602
$ws.Range("R20").Value = 6934793.265219755
Set-TextValue $ws.Range("Y20") "2023-08-18"
Set-TextValue $ws.Range("AA20") "2023-08-18"
$ws.Range("AC20").ClearContents()

# ---- Row 21 ----
$ws.Range("A21").Value = 111560058
$ws.Range("B21").Value = 78578
$ws.Range("E21").Value = 6458
Set-TextValue $ws.Range("F21") "Lunglav"
Set-TextValue $ws.Range("G21") "Lobaria pulmonaria"
Set-TextValue $ws.Range("H21") "(L.) Hoffm."
$ws.Range("Q21").Value = 523906.9737172622
$ws.Range("R21").Value = 6934619.326478666

# ---- Row 22 ----
$ws.Range("A22").Value = 111559701
$ws.Range("B22").Value = 12249
Set-TextValue $ws.Range("D22") "EN"
$ws.Range("E22").Value = 101283
Set-TextValue $ws.Range("F22") "Djupsvart brunbagge"
Set-TextValue $ws.Range("G22") "Melandrya dubia"
Set-TextValue $ws.Range("H22") "(Schaller, 1783)"
$ws.Range("Q22").Value = 523950.9321204902
$ws.Range("R22").Value = 6934675.944620069
Set-TextValue $ws.Range("Y22") "2023-08-17"
Set-TextValue $ws.Range("AA22") "2023-08-17"
Set-TextValue $ws.Range("AC22") "Kläckhål med svartfärgade larvgångar på björkhögstubbe med levande fnösktickor. Naturskog norr om Vattensjöarna"
